# Atualização de bases das ligas, do dia: 18-02-2024 às 22:54
#
# The underlying data rows were re-ordered (match records re-sorted);
# for each affected pair/group of rows, every data column (B..AC) moves
# to a different row while the leading sequential id column (A) stays
# put. Implemented as captured-array swaps/rotations over the full
# B:AC row ranges so every value (numbers, text, shared strings) moves
# together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 27 and 28: simple swap -------------------------------------
$row27 = $ws.Range("B27:AC27").Value()
$row28 = $ws.Range("B28:AC28").Value()

$ws.Range("B27:AC27").Value = $row28
$ws.Range("B28:AC28").Value = $row27

# --- Rows 175, 176, 178, 179: 4-way rotation --------------------------
# new175 = old176, new176 = old178, new178 = old179, new179 = old175
$row175 = $ws.Range("B175:AC175").Value()
$row176 = $ws.Range("B176:AC176").Value()
$row178 = $ws.Range("B178:AC178").Value()
$row179 = $ws.Range("B179:AC179").Value()

$ws.Range("B175:AC175").Value = $row176
$ws.Range("B176:AC176").Value = $row178
$ws.Range("B178:AC178").Value = $row179
$ws.Range("B179:AC179").Value = $row175
